# Append the latest scraped price-history row to the tracking sheet.
# Sheet columns: A=Date, B=Price, C=Discount, D=Incredible
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append right after the last used row so this keeps working as the
# history grows with every scrape.
$newRow = $ws.UsedRange.Rows.Count + 1

$targetRange = $ws.Range("A" + $newRow + ":D" + $newRow)

# "2026-02-07", "1192500", "10" and "1" look like a date/numbers to Excel's
# auto-detection, but every column in this sheet is stored as plain text
# (shared strings), matching all the existing history rows. Force text
# entry, then drop the temporary cell formatting back to the sheet's
# normal style so the new row stays visually identical to the others.
$targetRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "1192500"
$ws.Cells.Item($newRow, 3).Value = "10"
$ws.Cells.Item($newRow, 4).Value = "1"

$targetRange.Style = "Normal"
